$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting/styles from the last existing row (451) down through the new rows (452-463)
# so the new rows inherit the same cell styles (date format in col A, text style elsewhere).
$src = $ws.Range("A451:AD451")
$dst = $ws.Range("A452:AD463")
$src.Copy($dst)

# Row 452
$ws.Range("A452").Value = 45200.940594467596
$ws.Range("B452").Value = 'scw0922@naver.com'
$ws.Range("C452").Value = '간호학과'
$ws.Range("D452").Value = 20236256
$ws.Range("E452").Value = '신채원'
$ws.Range("F452").Value = 3
$ws.Range("G452").Value = '2. 시세 정보는 약 1개월 간격으로 갱신된다.'
$ws.Range("H452").Value = '2. 인공적인 향기가 인체에 해롭지 않을까요?'
$ws.Range("I452").Value = '4. 6630번 버스를 타면 한마음예식장에 갈 수 있다.'
$ws.Range("J452").Value = '3. 거실 바닥을 자주 물걸레로 닦는다'
$ws.Range("K452").Value = '3. 음주 운전이 의심될 경우 경찰관은 바로 운전자에게 혈액 채취를 명할 수 있군.'
$ws.Range("L452").Value = '2. 친구를 만나서 가까운 산에 오른다.'
$ws.Range("M452").Value = '3. 허위로 신고하면 10만 원의 과태료를 물게 된다.'
$ws.Range("N452").Value = '4. 8,000 원'
$ws.Range("O452").Value = '3. 기침 감기에 종합 감기약을 먹으면 기침약을 먹은 것과 효과가 같다.'
$ws.Range("P452").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q452").Value = '3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요.'
$ws.Range("R452").Value = '3. 홍길동 씨가 보내려는 돈은 30,500 원이다.'
$ws.Range("S452").Value = '2. 보증 기간 동안에는 건전지를 무상으로 제공한다.'
$ws.Range("T452").Value = '2. 동남쪽에서부터 꽃이 피기 시작한다.'
$ws.Range("U452").Value = '2. 벽지를 구입한 고객에게는 대걸레를 준다.'
$ws.Range("V452").Value = '2. 오전 6시'
$ws.Range("W452").Value = '1. 내 전공이 화학이니 지원 가능하겠군.'
$ws.Range("X452").Value = '3. 상점 주소'
$ws.Range("Y452").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z452").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA452").Value = '2. 교육은 특정 요일에 실시된다.'
$ws.Range("AB452").Value = '3. 공공장소에서는 전화 예절을 지켜야 한다.'
$ws.Range("AC452").Value = '4. 16권의 만화 ''토지''가 한꺼번에 출간되었다.'
$ws.Range("AD452").Value = '4. 야구 중계는 오후 2시 25분에 시작한다.'
$ws.Rows.Item(452).RowHeight = 15.75

# Row 453
$ws.Range("A453").Value = 45200.943949502311
$ws.Range("B453").Value = 'harin3040@naver.com'
$ws.Range("C453").Value = '심리학과'
$ws.Range("D453").Value = 20232113
$ws.Range("E453").Value = '김현진'
$ws.Range("F453").Value = 3
$ws.Range("G453").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H453").Value = '4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?'
$ws.Range("I453").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J453").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K453").Value = '2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군.'
$ws.Range("L453").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M453").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N453").Value = '3. 7,000 원'
$ws.Range("O453").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P453").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q453").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R453").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S453").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T453").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U453").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V453").Value = '4. 오후3시'
$ws.Range("W453").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X453").Value = '1. 상품 가격'
$ws.Range("Y453").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z453").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA453").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB453").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC453").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD453").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(453).RowHeight = 15.75

# Row 454
$ws.Range("A454").Value = 45200.946542951388
$ws.Range("B454").Value = 'shanesun0923@gmail.com'
$ws.Range("C454").Value = '간호학과'
$ws.Range("D454").Value = 20236253
$ws.Range("E454").Value = '선세인'
$ws.Range("F454").Value = 3
$ws.Range("G454").Value = '3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다.'
$ws.Range("H454").Value = '4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?'
$ws.Range("I454").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J454").Value = '3. 거실 바닥을 자주 물걸레로 닦는다'
$ws.Range("K454").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L454").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M454").Value = '1. 이 서식은 국내 전입신고 시에만 사용할 수 있다.'
$ws.Range("N454").Value = '3. 7,000 원'
$ws.Range("O454").Value = '4. 남은 약은 반드시 냉장고에 보관해야 한다.'
$ws.Range("P454").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q454").Value = '4. 확인증을 잃어버렸는데, 다시 발급해 주겠지.'
$ws.Range("R454").Value = '3. 홍길동 씨가 보내려는 돈은 30,500 원이다.'
$ws.Range("S454").Value = '2. 보증 기간 동안에는 건전지를 무상으로 제공한다.'
$ws.Range("T454").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U454").Value = '4. 조기 품절 시에는 사은품이 다른 물품으로 대체될 수 있다.'
$ws.Range("V454").Value = '3. 오후 6시'
$ws.Range("W454").Value = '4. 일주일에 이틀은 쉴 수 있겠군.'
$ws.Range("X454").Value = '3. 상점 주소'
$ws.Range("Y454").Value = '1. 뜻풀이 ''1'''
$ws.Range("Z454").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA454").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB454").Value = '4. 공공장소에서는 떠들지 말아야 한다.'
$ws.Range("AC454").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD454").Value = '4. 야구 중계는 오후 2시 25분에 시작한다.'
$ws.Rows.Item(454).RowHeight = 15.75

# Row 455
$ws.Range("A455").Value = 45200.95309476852
$ws.Range("B455").Value = 'sung93716@gmail.com'
$ws.Range("C455").Value = '데이터사이언스학부'
$ws.Range("D455").Value = 20233261
$ws.Range("E455").Value = '한예림'
$ws.Range("F455").Value = 3
$ws.Range("G455").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H455").Value = '1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?'
$ws.Range("I455").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J455").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K455").Value = '2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군.'
$ws.Range("L455").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M455").Value = '3. 허위로 신고하면 10만 원의 과태료를 물게 된다.'
$ws.Range("N455").Value = '2. 6,000 원'
$ws.Range("O455").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P455").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q455").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R455").Value = '4. 홍길동 씨는 세계은행에서 송금을 하고 있다.'
$ws.Range("S455").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T455").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U455").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V455").Value = '4. 오후3시'
$ws.Range("W455").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X455").Value = '1. 상품 가격'
$ws.Range("Y455").Value = '2. 뜻풀이 ''2'''
$ws.Range("Z455").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA455").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB455").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC455").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD455").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(455).RowHeight = 15.75

# Row 456
$ws.Range("A456").Value = 45200.968275844905
$ws.Range("B456").Value = 'dksdksqh1018@naver.com'
$ws.Range("C456").Value = '미디어스쿨'
$ws.Range("D456").Value = 20232549
$ws.Range("E456").Value = '안보민'
$ws.Range("F456").Value = 3
$ws.Range("G456").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H456").Value = '4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?'
$ws.Range("I456").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J456").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K456").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L456").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M456").Value = '3. 허위로 신고하면 10만 원의 과태료를 물게 된다.'
$ws.Range("N456").Value = '2. 6,000 원'
$ws.Range("O456").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P456").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q456").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R456").Value = '4. 홍길동 씨는 세계은행에서 송금을 하고 있다.'
$ws.Range("S456").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T456").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U456").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V456").Value = '4. 오후3시'
$ws.Range("W456").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X456").Value = '1. 상품 가격'
$ws.Range("Y456").Value = '1. 뜻풀이 ''1'''
$ws.Range("Z456").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA456").Value = '2. 교육은 특정 요일에 실시된다.'
$ws.Range("AB456").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC456").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD456").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(456).RowHeight = 15.75

# Row 457
$ws.Range("A457").Value = 45200.971810937495
$ws.Range("B457").Value = 'gaejisub@gmail.com'
$ws.Range("C457").Value = '콘텐츠it'
$ws.Range("D457").Value = 20225169
$ws.Range("E457").Value = '배승유'
$ws.Range("F457").Value = 3
$ws.Range("G457").Value = '3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다.'
$ws.Range("H457").Value = '1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?'
$ws.Range("I457").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J457").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K457").Value = '3. 음주 운전이 의심될 경우 경찰관은 바로 운전자에게 혈액 채취를 명할 수 있군.'
$ws.Range("L457").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M457").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N457").Value = '2. 6,000 원'
$ws.Range("O457").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P457").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q457").Value = '3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요.'
$ws.Range("R457").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S457").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T457").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U457").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V457").Value = '4. 오후3시'
$ws.Range("W457").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X457").Value = '1. 상품 가격'
$ws.Range("Y457").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z457").Value = '3. 여우비, 소나무향기'
$ws.Range("AA457").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB457").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC457").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD457").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(457).RowHeight = 15.75

# Row 458
$ws.Range("A458").Value = 45200.981516238426
$ws.Range("B458").Value = 'gustj1654@naver.com'
$ws.Range("C458").Value = '심리학과'
$ws.Range("D458").Value = 20232137
$ws.Range("E458").Value = '조현서'
$ws.Range("F458").Value = 3
$ws.Range("G458").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H458").Value = '1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?'
$ws.Range("I458").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J458").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K458").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L458").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M458").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N458").Value = '2. 6,000 원'
$ws.Range("O458").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P458").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q458").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R458").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S458").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T458").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U458").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V458").Value = '4. 오후3시'
$ws.Range("W458").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X458").Value = '1. 상품 가격'
$ws.Range("Y458").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z458").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA458").Value = '1. 이 프로그램은 노인들만을 위한 것이다.'
$ws.Range("AB458").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC458").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD458").Value = '3. ''TV는 사랑을 싣고''는 다시 보기를 제공하지 않는다.'
$ws.Rows.Item(458).RowHeight = 15.75

# Row 459
$ws.Range("A459").Value = 45200.986926377314
$ws.Range("B459").Value = 'yejin4259@naver.com'
$ws.Range("C459").Value = '언어청각학부'
$ws.Range("D459").Value = 20233951
$ws.Range("E459").Value = '이예진'
$ws.Range("F459").Value = 2
$ws.Range("G459").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H459").Value = '4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?'
$ws.Range("I459").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J459").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K459").Value = '2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군.'
$ws.Range("L459").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M459").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N459").Value = '1. 5,000 원'
$ws.Range("O459").Value = '4. 남은 약은 반드시 냉장고에 보관해야 한다.'
$ws.Range("P459").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q459").Value = '3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요.'
$ws.Range("R459").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S459").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T459").Value = '2. 동남쪽에서부터 꽃이 피기 시작한다.'
$ws.Range("U459").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V459").Value = '1. 오전 3시'
$ws.Range("W459").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X459").Value = '1. 상품 가격'
$ws.Range("Y459").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z459").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA459").Value = '2. 교육은 특정 요일에 실시된다.'
$ws.Range("AB459").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC459").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD459").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(459).RowHeight = 15.75

# Row 460
$ws.Range("A460").Value = 45200.990127766199
$ws.Range("B460").Value = 'ghys1837@naver.com'
$ws.Range("C460").Value = '언어청각학부 청각학 전공'
$ws.Range("D460").Value = 20213939
$ws.Range("E460").Value = '안영서'
$ws.Range("F460").Value = 3
$ws.Range("G460").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H460").Value = '2. 인공적인 향기가 인체에 해롭지 않을까요?'
$ws.Range("I460").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J460").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K460").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L460").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M460").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N460").Value = '2. 6,000 원'
$ws.Range("O460").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P460").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q460").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R460").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S460").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T460").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U460").Value = '2. 벽지를 구입한 고객에게는 대걸레를 준다.'
$ws.Range("V460").Value = '1. 오전 3시'
$ws.Range("W460").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X460").Value = '1. 상품 가격'
$ws.Range("Y460").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z460").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA460").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB460").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC460").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD460").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(460).RowHeight = 15.75

# Row 461
$ws.Range("A461").Value = 45200.998866030088
$ws.Range("B461").Value = 'tjdus3641@gmail.com'
$ws.Range("C461").Value = '간호학과'
$ws.Range("D461").Value = 20226283
$ws.Range("E461").Value = '장서연'
$ws.Range("F461").Value = 3
$ws.Range("G461").Value = '3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다.'
$ws.Range("H461").Value = '4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?'
$ws.Range("I461").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J461").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K461").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L461").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M461").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N461").Value = '2. 6,000 원'
$ws.Range("O461").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P461").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q461").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R461").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S461").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T461").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U461").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V461").Value = '4. 오후3시'
$ws.Range("W461").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X461").Value = '1. 상품 가격'
$ws.Range("Y461").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z461").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA461").Value = '2. 교육은 특정 요일에 실시된다.'
$ws.Range("AB461").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC461").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD461").Value = '4. 야구 중계는 오후 2시 25분에 시작한다.'
$ws.Rows.Item(461).RowHeight = 15.75

# Row 462
$ws.Range("A462").Value = 45201.002475104164
$ws.Range("B462").Value = 'rkqls3333@gmail.com'
$ws.Range("C462").Value = '간호학과'
$ws.Range("D462").Value = 20236205
$ws.Range("E462").Value = '권가빈'
$ws.Range("F462").Value = 3
$ws.Range("G462").Value = '3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다.'
$ws.Range("H462").Value = '1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?'
$ws.Range("I462").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J462").Value = '3. 거실 바닥을 자주 물걸레로 닦는다'
$ws.Range("K462").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L462").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M462").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N462").Value = '2. 6,000 원'
$ws.Range("O462").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P462").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q462").Value = '2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?'
$ws.Range("R462").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S462").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T462").Value = '3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다.'
$ws.Range("U462").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V462").Value = '1. 오전 3시'
$ws.Range("W462").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X462").Value = '1. 상품 가격'
$ws.Range("Y462").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z462").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA462").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB462").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC462").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD462").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(462).RowHeight = 15.75

# Row 463
$ws.Range("A463").Value = 45201.009580462967
$ws.Range("B463").Value = 'nanadiana222@naver.com'
$ws.Range("C463").Value = '빅데이터과'
$ws.Range("D463").Value = 20225261
$ws.Range("E463").Value = '조희진'
$ws.Range("F463").Value = 3
$ws.Range("G463").Value = '4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다.'
$ws.Range("H463").Value = '1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?'
$ws.Range("I463").Value = '3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다.'
$ws.Range("J463").Value = '2. 채소 등 농산물은 익혀 먹는다'
$ws.Range("K463").Value = '4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군.'
$ws.Range("L463").Value = '1. 학교 가는 아이에게 우산을 챙겨 준다.'
$ws.Range("M463").Value = '4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다.'
$ws.Range("N463").Value = '2. 6,000 원'
$ws.Range("O463").Value = '2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다.'
$ws.Range("P463").Value = '1. 이 공지 사항은 ''봉투나라''의 관리자가 작성한 것이다.'
$ws.Range("Q463").Value = '3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요.'
$ws.Range("R463").Value = '2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다.'
$ws.Range("S463").Value = '4. 수리할 수 없는 제품은 새것으로 교환해 준다.'
$ws.Range("T463").Value = '2. 동남쪽에서부터 꽃이 피기 시작한다.'
$ws.Range("U463").Value = '1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다.'
$ws.Range("V463").Value = '4. 오후3시'
$ws.Range("W463").Value = '2. 남녀를 차별 대우하다니 한심하네.'
$ws.Range("X463").Value = '1. 상품 가격'
$ws.Range("Y463").Value = '3. 뜻풀이 ''3'''
$ws.Range("Z463").Value = '2. 여우비, 장난감박물관'
$ws.Range("AA463").Value = '4. 5만원을 내면 모든 강의를 들을 수 있다.'
$ws.Range("AB463").Value = '1. 출산율을 높여야 한다.'
$ws.Range("AC463").Value = '3. ''토지''에는 700여 명의 인물이 등장한다.'
$ws.Range("AD463").Value = '2. 낮 12시 뉴스는 자막 방송을 하지 않는다.'
$ws.Rows.Item(463).RowHeight = 15.75

# Update selection to the last entered cell, matching the final author state
$ws.Range("AD463").Select()
